$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.467.28'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').Value = '2.163.23'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.18'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.10%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0866'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.71%  '
$ws.Range('D13').Value = '2.483.08'
$ws.Range('E13').Value = '  +3.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.817'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').Value = '2.165.07'
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').Value = '39.456.10'
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.47%  '
$ws.Range('E25').Value = '  -2.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  +3.00%  '
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.07%  '
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.73%  '
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.47%  '
$ws.Range('E36').Value = '  +2.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.43'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.62'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '104.68'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.21%  '
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '1.538.19'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('E44').Value = '  +6.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0934'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.83%  '
$ws.Range('E47').Value = '  +7.43%  '
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').Value = '2.366.97'
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('E51').Value = '  +0.33%  '
